$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.008.82'
$ws.Range("E2").Value = '  -0.62%  '
$ws.Range("D3").Value = '2.611.24'
$ws.Range("E3").Value = '  -1.21%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''590.52'
$ws.Range("E5").Value = '  -1.37%  '
$ws.Range("D6").Value = '''165.30'
$ws.Range("E6").Value = '  -0.52%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '''0.531'
$ws.Range("E8").Value = '  -2.13%  '
$ws.Range("D9").Value = '2.610.60'
$ws.Range("E9").Value = '  -1.17%  '
$ws.Range("D10").Value = '''0.137'
$ws.Range("E10").Value = '  -5.56%  '
$ws.Range("E11").Value = '  +1.55%  '
$ws.Range("D12").Value = '''0.362'
$ws.Range("E12").Value = '  -0.55%  '
$ws.Range("D13").Value = '''5.19'
$ws.Range("E13").Value = '  -0.68%  '
$ws.Range("D14").Value = '''27.24'
$ws.Range("E14").Value = '  -2.57%  '
$ws.Range("D15").Value = '3.088.95'
$ws.Range("E15").Value = '  -1.01%  '
$ws.Range("D16").Value = '''0.0000179'
$ws.Range("E16").Value = '  -2.97%  '
$ws.Range("D17").Value = '67.120.60'
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("D18").Value = '2.617.65'
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("D19").Value = '''11.76'
$ws.Range("E19").Value = '  -0.84%  '
$ws.Range("D20").Value = '''7.81'
$ws.Range("E20").Value = '  -0.68%  '
$ws.Range("D21").Value = '''353.54'
$ws.Range("E21").Value = '  -2.95%  '
$ws.Range("D22").Value = '''4.27'
$ws.Range("E22").Value = '  -2.92%  '
$ws.Range("D23").Value = '''4.61'
$ws.Range("E23").Value = '  -3.57%  '
$ws.Range("D24").Value = '''10.49'
$ws.Range("E24").Value = '  -4.95%  '
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("E26").Value = '  -4.44%  '
$ws.Range("D27").Value = '''68.92'
$ws.Range("E27").Value = '  -2.60%  '
$ws.Range("D28").Value = '2.750.59'
$ws.Range("E28").Value = '  -0.95%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("D30").Value = '0.0₃0994'
$ws.Range("E30").Value = '  -3.14%  '
$ws.Range("D31").Value = '''538.82'
$ws.Range("E31").Value = '  -2.96%  '
$ws.Range("D32").Value = '''7.88'
$ws.Range("E32").Value = '  -1.81%  '
$ws.Range("E33").Value = '  -3.87%  '
$ws.Range("D34").Value = '''1.87'
$ws.Range("E34").Value = '  -2.45%  '
$ws.Range("D35").Value = '''0.134'
$ws.Range("E35").Value = '  +1.15%  '
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("D37").Value = '''1.48'
$ws.Range("E37").Value = '  -3.81%  '
$ws.Range("D38").Value = '''157.02'
$ws.Range("E38").Value = '  -0.50%  '
$ws.Range("D39").Value = '''18.89'
$ws.Range("E39").Value = '  -2.68%  '
$ws.Range("E40").Value = '  -2.26%  '
$ws.Range("E41").Value = '  +1.64%  '
$ws.Range("E42").Value = '  -1.46%  '
$ws.Range("D43").Value = '''5.12'
$ws.Range("E43").Value = '  -2.77%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = '''2.40'
$ws.Range("E45").Value = '  -5.11%  '
$ws.Range("D46").Value = '0.0₆0300'
$ws.Range("E46").Value = '  -1.13%  '
$ws.Range("D47").Value = '''151.08'
$ws.Range("E47").Value = '  -1.90%  '
$ws.Range("D48").Value = '''0.573'
$ws.Range("E48").Value = '  -3.69%  '
$ws.Range("D49").Value = '''3.76'
$ws.Range("E49").Value = '  -2.89%  '
$ws.Range("D50").Value = '''1.69'
$ws.Range("E50").Value = '  -2.32%  '
$ws.Range("E51").Value = '  -1.19%  '
